$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Half_Rocker_Width (mm) value from 7.5 to 10
$ws.Range("B17").Value = 10

# Delete the entire "Suspension_Rod_Rext (mm)" row (row 21), shifting rows below upward
$ws.Rows("21").Delete()

# Update the active selection to match the final state
$ws.Range("B18").Select()
